# "subimos el último SPA" - update last-week target units (R) and resulting
# consumption trend (T) / final order (U) for a batch of articles, and hide
# the now-fully-covered row 33 (SUBSTRATO PLANTAS VERDES 10L), then refresh
# the Total_Unidades summary metric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new value for column R ("uds. Objetivo semana pasada")
$rUpdates = @{
    6  = 5
    7  = 2
    8  = 6
    9  = 7
    10 = 14
    11 = 1
    12 = 1
    13 = 2
    14 = 6
    24 = 3
    26 = 2
    29 = 5
    33 = 1
    36 = 13
    37 = 7
    38 = 6
    39 = 2
    40 = 3
    41 = 6
    42 = 4
    43 = 1
}

foreach ($row in $rUpdates.Keys) {
    $ws.Range("R$row").Value = $rUpdates[$row]
}

# Row -> new value for column T ("Tendencia Consumo", 20=19-18, floored at 0)
$tUpdates = @{
    6  = 0
    9  = 0
    10 = 0
    12 = 1
    13 = 1
    14 = 0
    26 = 0
    29 = 0
    33 = 7
    36 = 0
    38 = 0
    41 = 10
    42 = 0
}

foreach ($row in $tUpdates.Keys) {
    $ws.Range("T$row").Value = $tUpdates[$row]
}

# Row 33 (SUBSTRATO PLANTAS VERDES 10L) is now fully stocked - its "Pedido
# Final" order quantity drops to 0 and the row gets hidden from view.
$ws.Range("U33").Value = 0
$ws.Rows(33).EntireRow.Hidden = $true

# Refresh the summary metric Total_Unidades (= SUM of column U across the
# article rows), which moves from 41 to 40 because of the U33 change.
$ws.Range("C47").Value = 40
